$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.072892327087841
$ws.Range("D2").Value = 1.071870215113769
$ws.Range("E2").Value = 1.076205767245856
$ws.Range("F2").Value = 1.084854449401715
$ws.Range("I2").Value = 1.052804418647452
$ws.Range("J2").Value = 1.077809583936983
$ws.Range("K2").Value = 1.074566238892568
$ws.Range("L2").Value = 1.078890311757841
$ws.Range("M2").Value = 1.087516390636248

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.074397094034832
$ws.Range("D3").Value = 1.073022815664558
$ws.Range("E3").Value = 1.077513318269282
$ws.Range("F3").Value = 1.086160823543457
$ws.Range("I3").Value = 1.053211707738529
$ws.Range("J3").Value = 1.078970231807105
$ws.Range("K3").Value = 1.075534646404976
$ws.Range("L3").Value = 1.080014124295421
$ws.Range("M3").Value = 1.088640674898526

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.075369739712432
$ws.Range("D4").Value = 1.073767481534873
$ws.Range("E4").Value = 1.078358574765837
$ws.Range("F4").Value = 1.087005272722475
$ws.Range("I4").Value = 1.053473160082334
$ws.Range("J4").Value = 1.079719760781945
$ws.Range("K4").Value = 1.076159544114203
$ws.Range("L4").Value = 1.080739936370978
$ws.Range("M4").Value = 1.089366743467913

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.075778397506665
$ws.Range("D5").Value = 1.074080269215919
$ws.Range("E5").Value = 1.07871372982578
$ws.Range("F5").Value = 1.087360077275099
$ws.Range("I5").Value = 1.053582575919267
$ws.Range("J5").Value = 1.080034511184373
$ws.Range("K5").Value = 1.076421841486718
$ws.Range("L5").Value = 1.081044744174707
$ws.Range("M5").Value = 1.089671647303483

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.075846998914591
$ws.Range("D6").Value = 1.074132771934592
$ws.Range("E6").Value = 1.078773350872029
$ws.Range("F6").Value = 1.087419638813768
$ws.Range("I6").Value = 1.053600918127457
$ws.Range("J6").Value = 1.080087338678589
$ws.Range("K6").Value = 1.076465858457821
$ws.Range("L6").Value = 1.08109590389616
$ws.Range("M6").Value = 1.089722822453485

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.075375201162691
$ws.Range("D7").Value = 1.073771662074601
$ws.Range("E7").Value = 1.07836332110892
$ws.Range("F7").Value = 1.087010014425908
$ws.Range("I7").Value = 1.053474624059493
$ws.Range("J7").Value = 1.079723967867994
$ws.Range("K7").Value = 1.076163050549442
$ws.Range("L7").Value = 1.080744010493901
$ws.Range("M7").Value = 1.089370818920541

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.0734010891039
$ws.Range("D8").Value = 1.072259981083573
$ws.Range("E8").Value = 1.076647832115718
$ws.Range("F8").Value = 1.085296125933075
$ws.Range("I8").Value = 1.052942497958399
$ws.Range("J8").Value = 1.078202141157856
$ws.Range("K8").Value = 1.07489387653637
$ws.Range("L8").Value = 1.079270395705258
$ws.Range("M8").Value = 1.087896644020755

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.069914208633012
$ws.Range("D9").Value = 1.06958726764389
$ws.Range("E9").Value = 1.073618457584905
$ws.Range("F9").Value = 1.082269227434136
$ws.Range("I9").Value = 1.051988726531689
$ws.Range("J9").Value = 1.075508886204515
$ws.Range("K9").Value = 1.072644036389282
$ws.Range("L9").Value = 1.076663012365921
$ws.Range("M9").Value = 1.085287906042564

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.067583680642252
$ws.Range("D10").Value = 1.067799196899214
$ws.Range("E10").Value = 1.071594231933961
$ws.Range("F10").Value = 1.080246429856689
$ws.Range("I10").Value = 1.051341938204248
$ws.Range("J10").Value = 1.073705290956814
$ws.Range("K10").Value = 1.071134899762848
$ws.Range("L10").Value = 1.074917297792348
$ws.Range("M10").Value = 1.08354104888135

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.066573037837228
$ws.Range("D11").Value = 1.067023400473677
$ws.Range("E11").Value = 1.070716552022095
$ws.Range("F11").Value = 1.079369317412503
$ws.Range("I11").Value = 1.051059250038795
$ws.Range("J11").Value = 1.072922329509029
$ws.Range("K11").Value = 1.070479180912981
$ws.Range("L11").Value = 1.074159553250237
$ws.Range("M11").Value = 1.082782754423863

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.066197405645003
$ws.Range("D12").Value = 1.066734997278332
$ws.Range("E12").Value = 1.070390359985739
$ws.Range("F12").Value = 1.079043328669163
$ws.Range("I12").Value = 1.050953850385978
$ws.Range("J12").Value = 1.072631197638986
$ws.Range("K12").Value = 1.07023527457366
$ws.Range("L12").Value = 1.07387781155888
$ws.Range("M12").Value = 1.082500800316376

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.066277990760315
$ws.Range("D13").Value = 1.066796871586975
$ws.Range("E13").Value = 1.070460337643629
$ws.Range("F13").Value = 1.079113263056426
$ws.Range("I13").Value = 1.050976476958253
$ws.Range("J13").Value = 1.072693660313448
$ws.Range("K13").Value = 1.070287608916304
$ws.Range("L13").Value = 1.073938258939491
$ws.Range("M13").Value = 1.08256129362775

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.06654199279101
$ws.Range("D14").Value = 1.066999565856871
$ws.Range("E14").Value = 1.070689592646247
$ws.Range("F14").Value = 1.079342374992802
$ws.Range("I14").Value = 1.051050545776584
$ws.Range("J14").Value = 1.072898270715618
$ws.Range("K14").Value = 1.070459026566695
$ws.Range("L14").Value = 1.074136270170664
$ws.Range("M14").Value = 1.082759453951448

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.066704621854236
$ws.Range("D15").Value = 1.067124420817523
$ws.Range("E15").Value = 1.070830819843289
$ws.Range("F15").Value = 1.07948351305423
$ws.Range("I15").Value = 1.051096129384712
$ws.Range("J15").Value = 1.073024297333009
$ws.Range("K15").Value = 1.070564597024593
$ws.Range("L15").Value = 1.074258233944569
$ws.Range("M15").Value = 1.082881508509844

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.067650721225466
$ws.Range("D16").Value = 1.067850650871158
$ws.Range("E16").Value = 1.071652455367562
$ws.Range("F16").Value = 1.080304614570397
$ws.Range("I16").Value = 1.051360643788426
$ws.Range("J16").Value = 1.073757211055133
$ws.Range("K16").Value = 1.071178369799198
$ws.Range("L16").Value = 1.074967547604567
$ws.Range("M16").Value = 1.083591334038889

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.068243774941521
$ws.Range("D17").Value = 1.068305777306999
$ws.Range("E17").Value = 1.072167526444862
$ws.Range("F17").Value = 1.080819337162912
$ws.Range("I17").Value = 1.051525862349306
$ws.Range("J17").Value = 1.074216410905512
$ws.Range("K17").Value = 1.071562766806246
$ws.Range("L17").Value = 1.075411985274198
$ws.Range("M17").Value = 1.084036078116318

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.068589548090895
$ws.Range("D18").Value = 1.068571095644893
$ws.Range("E18").Value = 1.07246784553789
$ws.Range("F18").Value = 1.081119448069274
$ws.Range("I18").Value = 1.051621978471062
$ws.Range("J18").Value = 1.074484062522349
$ws.Range("K18").Value = 1.071786762041963
$ws.Range("L18").Value = 1.075671041215486
$ws.Range("M18").Value = 1.084295307448207

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.068707423430479
$ws.Range("D19").Value = 1.068661537182835
$ws.Range("E19").Value = 1.072570227570577
$ws.Range("F19").Value = 1.081221758276154
$ws.Range("I19").Value = 1.051654708731043
$ws.Range("J19").Value = 1.07457529243209
$ws.Range("K19").Value = 1.07186310195772
$ws.Range("L19").Value = 1.075759342731633
$ws.Range("M19").Value = 1.084383667177568

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.068180160949973
$ws.Range("D20").Value = 1.068256962017677
$ws.Range("E20").Value = 1.072112275939056
$ws.Range("F20").Value = 1.080764124554131
$ws.Range("I20").Value = 1.051508162154744
$ws.Range("J20").Value = 1.074167162996764
$ws.Range("K20").Value = 1.071521547128487
$ws.Range("L20").Value = 1.075364319668844
$ws.Range("M20").Value = 1.08398838018683

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.066464257307851
$ws.Range("D21").Value = 1.066939884052796
$ws.Range("E21").Value = 1.07062208785074
$ws.Range("F21").Value = 1.079274912532233
$ws.Range("I21").Value = 1.051028745325196
$ws.Range("J21").Value = 1.072838026521174
$ws.Range("K21").Value = 1.07040855785213
$ws.Range("L21").Value = 1.074077968616285
$ws.Range("M21").Value = 1.082701108717021

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.065384040209875
$ws.Range("D22").Value = 1.066110406174395
$ws.Range("E22").Value = 1.069684087848483
$ws.Range("F22").Value = 1.078337482965491
$ws.Range("I22").Value = 1.050725020671367
$ws.Range("J22").Value = 1.072000577428419
$ws.Range("K22").Value = 1.069706789526602
$ws.Range("L22").Value = 1.073267556530596
$ws.Range("M22").Value = 1.081890070823554

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.065956815287211
$ws.Range("D23").Value = 1.066550260451003
$ws.Range("E23").Value = 1.070181441791566
$ws.Range("F23").Value = 1.078834538557879
$ws.Range("I23").Value = 1.050886249312861
$ws.Range("J23").Value = 1.072444694568654
$ws.Range("K23").Value = 1.070079000196045
$ws.Range("L23").Value = 1.073697327828436
$ws.Range("M23").Value = 1.082320178294044

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.068208905823034
$ws.Range("D24").Value = 1.068279020007081
$ws.Range("E24").Value = 1.072137241614376
$ws.Range("F24").Value = 1.080789073120636
$ws.Range("I24").Value = 1.051516160892885
$ws.Range("J24").Value = 1.07418941660009
$ws.Range("K24").Value = 1.071540173197485
$ws.Range("L24").Value = 1.075385858251475
$ws.Range("M24").Value = 1.084009933391833

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.070816668567939
$ws.Range("D25").Value = 1.070279313222119
$ws.Range("E25").Value = 1.074402419940958
$ws.Range("F25").Value = 1.083052589558076
$ws.Range("I25").Value = 1.052237218986952
$ws.Range("J25").Value = 1.076206561943776
$ws.Range("K25").Value = 1.073227285605249
$ws.Range("L25").Value = 1.077338377933488
$ws.Range("M25").Value = 1.085963664588901

Write-Output "Applied 240 cell updates"